$wb = $excel.ActiveWorkbook

# --- Projects sheet: drop the three sample/test rows (rows 2-4), shifting
#     the sheet back down to just its header row. ---
$wsProjects = $wb.Worksheets.Item("Projects")
[void]$wsProjects.Activate()
[void]$wsProjects.Rows("2:4").Delete()
[void]$wsProjects.Range("F14").Select()

# --- Results / Impacts: selection moved to F14 as well. ---
$wsResults = $wb.Worksheets.Item("Results")
[void]$wsResults.Activate()
[void]$wsResults.Range("F14").Select()

$wsImpacts = $wb.Worksheets.Item("Impacts")
[void]$wsImpacts.Activate()
[void]$wsImpacts.Range("F14").Select()

# --- Rename "Contributions" -> "Donor Contribution". ---
$wsContrib = $wb.Worksheets.Item("Contributions")
$wsContrib.Name = "Donor Contribution"

# --- Make it the active/selected sheet with its new selection. ---
[void]$wsContrib.Activate()
[void]$wsContrib.Range("E15").Select()
